$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "POPULATION"
$ws.Range("C2").Value = "13.3 M"
$ws.Range("C3").Value = "15.4 M"
$ws.Range("C4").Value = "17.5 M"
$ws.Range("C5").Value = "19.6 M"
$ws.Range("C6").Value = "111.7 M"
$ws.Range("C7").Value = "113.8 M"
$ws.Range("C8").Value = "115.9 M"
$ws.Range("C9").Value = "117.10 M"
$ws.Range("C10").Value = "119.11 M"
$ws.Range("C11").Value = "121.12 M"

$ws.Range("C1").Select() | Out-Null
